$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "joint-datasets_after_task1"

# Update metric values for epochs 1-40 (rows 3-42)
$ws.Range("C3").Value2 = 1.015314094908846
$ws.Range("D3").Value2 = 0.4915474355220795
$ws.Range("E3").Value2 = 85.11
$ws.Range("I3").Value2 = 0.002111403516493738
$ws.Range("J3").Value2 = 85.11
$ws.Range("C4").Value2 = 0.6214573278849034
$ws.Range("D4").Value2 = 0.4267416298389435
$ws.Range("E4").Value2 = 87.015
$ws.Range("I4").Value2 = 0.001787694663065486
$ws.Range("J4").Value2 = 87.015
$ws.Range("C5").Value2 = 0.5542708210435654
$ws.Range("D5").Value2 = 0.3855285942554474
$ws.Range("E5").Value2 = 87.795
$ws.Range("I5").Value2 = 0.001671081167226657
$ws.Range("J5").Value2 = 87.795
$ws.Range("C6").Value2 = 0.5162223602873655
$ws.Range("D6").Value2 = 0.347023993730545
$ws.Range("E6").Value2 = 88.455
$ws.Range("I6").Value2 = 0.001572182883857749
$ws.Range("J6").Value2 = 88.455
$ws.Range("C7").Value2 = 0.491414361286641
$ws.Range("D7").Value2 = 0.3411555588245392
$ws.Range("E7").Value2 = 89.11
$ws.Range("I7").Value2 = 0.001487036141770659
$ws.Range("J7").Value2 = 89.11
$ws.Range("C8").Value2 = 0.4655047141152351
$ws.Range("D8").Value2 = 0.3380215167999268
$ws.Range("E8").Value2 = 89.67
$ws.Range("I8").Value2 = 0.001424678029294591
$ws.Range("J8").Value2 = 89.67
$ws.Range("C9").Value2 = 0.4544377939430421
$ws.Range("D9").Value2 = 0.3336684703826904
$ws.Range("E9").Value2 = 90.155
$ws.Range("I9").Value2 = 0.001384824806987308
$ws.Range("J9").Value2 = 90.155
$ws.Range("C10").Value2 = 0.4424802191269417
$ws.Range("D10").Value2 = 0.3207788169384003
$ws.Range("E10").Value2 = 90.255
$ws.Range("I10").Value2 = 0.001348726329897181
$ws.Range("J10").Value2 = 90.255
$ws.Range("C11").Value2 = 0.4306785971274559
$ws.Range("D11").Value2 = 0.3228722512722015
$ws.Range("E11").Value2 = 90.68000000000001
$ws.Range("I11").Value2 = 0.001314192454481963
$ws.Range("J11").Value2 = 90.68000000000001
$ws.Range("C12").Value2 = 0.4247346638339589
$ws.Range("D12").Value2 = 0.3134530186653137
$ws.Range("E12").Value2 = 90.84
$ws.Range("I12").Value2 = 0.001294319033416104
$ws.Range("J12").Value2 = 90.84
$ws.Range("C13").Value2 = 0.4159635641738051
$ws.Range("D13").Value2 = 0.316488653421402
$ws.Range("E13").Value2 = 90.98999999999999
$ws.Range("I13").Value2 = 0.001261701207174337
$ws.Range("J13").Value2 = 90.98999999999999
$ws.Range("C14").Value2 = 0.4097811999623485
$ws.Range("D14").Value2 = 0.3075658977031708
$ws.Range("E14").Value2 = 91.095
$ws.Range("I14").Value2 = 0.001247735364391701
$ws.Range("J14").Value2 = 91.095
$ws.Range("C15").Value2 = 0.4030830496191382
$ws.Range("D15").Value2 = 0.3101250529289246
$ws.Range("E15").Value2 = 91.395
$ws.Range("I15").Value2 = 0.001235743238031864
$ws.Range("J15").Value2 = 91.395
$ws.Range("C16").Value2 = 0.3999650908730464
$ws.Range("D16").Value2 = 0.300663560628891
$ws.Range("E16").Value2 = 91.575
$ws.Range("I16").Value2 = 0.001213396305545757
$ws.Range("J16").Value2 = 91.575
$ws.Range("C17").Value2 = 0.3939806901801607
$ws.Range("D17").Value2 = 0.2853456735610962
$ws.Range("E17").Value2 = 91.33
$ws.Range("I17").Value2 = 0.00119290139144141
$ws.Range("J17").Value2 = 91.33
$ws.Range("C18").Value2 = 0.3879827733628937
$ws.Range("D18").Value2 = 0.3020427823066711
$ws.Range("E18").Value2 = 91.55500000000001
$ws.Range("I18").Value2 = 0.001177788913589029
$ws.Range("J18").Value2 = 91.55500000000001
$ws.Range("C19").Value2 = 0.3878454062745647
$ws.Range("D19").Value2 = 0.2924025058746338
$ws.Range("E19").Value2 = 91.61
$ws.Range("I19").Value2 = 0.001159678466452169
$ws.Range("J19").Value2 = 91.61
$ws.Range("C20").Value2 = 0.3807221245337806
$ws.Range("D20").Value2 = 0.2904521226882935
$ws.Range("E20").Value2 = 91.88500000000001
$ws.Range("I20").Value2 = 0.001152909800200723
$ws.Range("J20").Value2 = 91.88500000000001
$ws.Range("C21").Value2 = 0.378092708019462
$ws.Range("D21").Value2 = 0.3071305453777313
$ws.Range("E21").Value2 = 92.095
$ws.Range("I21").Value2 = 0.001140092406867188
$ws.Range("J21").Value2 = 92.095
$ws.Range("C22").Value2 = 0.3723721927364601
$ws.Range("D22").Value2 = 0.2974560558795929
$ws.Range("E22").Value2 = 92.08499999999999
$ws.Range("I22").Value2 = 0.001126028405058605
$ws.Range("J22").Value2 = 92.08499999999999
$ws.Range("C23").Value2 = 0.3800110408977196
$ws.Range("D23").Value2 = 0.2988713383674622
$ws.Range("E23").Value2 = 91.69499999999999
$ws.Range("I23").Value2 = 0.001169747290216037
$ws.Range("J23").Value2 = 91.69499999999999
$ws.Range("C24").Value2 = 0.3737525993675143
$ws.Range("D24").Value2 = 0.2933715283870697
$ws.Range("E24").Value2 = 91.86
$ws.Range("I24").Value2 = 0.001161904219538701
$ws.Range("J24").Value2 = 91.86
$ws.Range("C25").Value2 = 0.3698755489516139
$ws.Range("D25").Value2 = 0.2941969633102417
$ws.Range("E25").Value2 = 91.715
$ws.Range("I25").Value2 = 0.001155274915140762
$ws.Range("J25").Value2 = 91.715
$ws.Range("C26").Value2 = 0.3675317036140343
$ws.Range("D26").Value2 = 0.292375385761261
$ws.Range("E26").Value2 = 91.995
$ws.Range("I26").Value2 = 0.001150498993581277
$ws.Range("J26").Value2 = 91.995
$ws.Range("C27").Value2 = 0.3672150623718765
$ws.Range("D27").Value2 = 0.287432849407196
$ws.Range("E27").Value2 = 91.90000000000001
$ws.Range("I27").Value2 = 0.001148534108557215
$ws.Range("J27").Value2 = 91.90000000000001
$ws.Range("C28").Value2 = 0.38178953697566
$ws.Range("D28").Value2 = 0.2936921119689941
$ws.Range("E28").Value2 = 91.675
$ws.Range("I28").Value2 = 0.001178276808322698
$ws.Range("J28").Value2 = 91.675
$ws.Range("C29").Value2 = 0.3792383054361518
$ws.Range("D29").Value2 = 0.2972906231880188
$ws.Range("E29").Value2 = 91.655
$ws.Range("I29").Value2 = 0.001175190215368639
$ws.Range("J29").Value2 = 91.655
$ws.Range("C30").Value2 = 0.3778147419078123
$ws.Range("D30").Value2 = 0.2987631857395172
$ws.Range("E30").Value2 = 91.655
$ws.Range("I30").Value2 = 0.001173762683779932
$ws.Range("J30").Value2 = 91.655
$ws.Range("C31").Value2 = 0.3739317507298044
$ws.Range("D31").Value2 = 0.2968397736549377
$ws.Range("E31").Value2 = 91.66
$ws.Range("I31").Value2 = 0.001168318200488284
$ws.Range("J31").Value2 = 91.66
$ws.Range("C32").Value2 = 0.3719542741725759
$ws.Range("D32").Value2 = 0.2977724671363831
$ws.Range("E32").Value2 = 91.68000000000001
$ws.Range("I32").Value2 = 0.001167451728068409
$ws.Range("J32").Value2 = 91.68000000000001
$ws.Range("C33").Value2 = 0.3824908246032782
$ws.Range("D33").Value2 = 0.2886776328086853
$ws.Range("E33").Value2 = 91.54000000000001
$ws.Range("I33").Value2 = 0.001185812481533503
$ws.Range("J33").Value2 = 91.54000000000001
$ws.Range("C34").Value2 = 0.3813699954290422
$ws.Range("D34").Value2 = 0.290938526391983
$ws.Range("E34").Value2 = 91.575
$ws.Range("I34").Value2 = 0.001183568340477359
$ws.Range("J34").Value2 = 91.575
$ws.Range("C35").Value2 = 0.3823275466345387
$ws.Range("D35").Value2 = 0.2921743988990784
$ws.Range("E35").Value2 = 91.62
$ws.Range("I35").Value2 = 0.001181084162379557
$ws.Range("J35").Value2 = 91.62
$ws.Range("C36").Value2 = 0.3818072689096995
$ws.Range("D36").Value2 = 0.2930486500263214
$ws.Range("E36").Value2 = 91.63500000000001
$ws.Range("I36").Value2 = 0.001180295584726264
$ws.Range("J36").Value2 = 91.63500000000001
$ws.Range("C37").Value2 = 0.3812154673276242
$ws.Range("D37").Value2 = 0.2937979400157928
$ws.Range("E37").Value2 = 91.66
$ws.Range("I37").Value2 = 0.001179288306734816
$ws.Range("J37").Value2 = 91.66
$ws.Range("C38").Value2 = 0.3875494549837654
$ws.Range("D38").Value2 = 0.2860315442085266
$ws.Range("E38").Value2 = 91.375
$ws.Range("I38").Value2 = 0.001190961829897424
$ws.Range("J38").Value2 = 91.375
$ws.Range("C39").Value2 = 0.3878712158071776
$ws.Range("D39").Value2 = 0.2867408096790314
$ws.Range("E39").Value2 = 91.405
$ws.Range("I39").Value2 = 0.001189245010349259
$ws.Range("J39").Value2 = 91.405
$ws.Range("C40").Value2 = 0.3857999317634086
$ws.Range("D40").Value2 = 0.2872774600982666
$ws.Range("E40").Value2 = 91.45999999999999
$ws.Range("I40").Value2 = 0.001188002932514064
$ws.Range("J40").Value2 = 91.45999999999999
$ws.Range("C41").Value2 = 0.3819398087929803
$ws.Range("D41").Value2 = 0.2877419888973236
$ws.Range("E41").Value2 = 91.485
$ws.Range("I41").Value2 = 0.00118685124819749
$ws.Range("J41").Value2 = 91.485
$ws.Range("C42").Value2 = 0.3842184259517364
$ws.Range("D42").Value2 = 0.2884418368339539
$ws.Range("E42").Value2 = 91.5
$ws.Range("I42").Value2 = 0.001185963048883423
$ws.Range("J42").Value2 = 91.5

# Remove the now-unused trailing epoch rows (41-53), which were rows 43-55
$ws.Range("A43:K55").EntireRow.Delete()
